# Rename the stats sheets to their human-readable display names.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("StandardStats")
$ws.Name = "Standard Stats"

$ws = $wb.Worksheets.Item("ShootingStats")
$ws.Name = "Shooting Stats"

$ws = $wb.Worksheets.Item("PassingStats")
$ws.Name = "Passing Stats"

$ws = $wb.Worksheets.Item("PassTypes")
$ws.Name = "Pass Types"

$ws = $wb.Worksheets.Item("GoalShotCreation")
$ws.Name = "Goal & Shot Creation"

$ws = $wb.Worksheets.Item("DefensiveActions")
$ws.Name = "Defensive Actions"

$ws = $wb.Worksheets.Item("PlayingTime")
$ws.Name = "Playing Time"

$ws = $wb.Worksheets.Item("MiscStats")
$ws.Name = "Miscellaneous Stats"

# Bump every player's "Age" column (format YY-DDD, years-days) forward by
# one day across every per-player stats sheet. "Possession" keeps its
# original tab name but still gets the same Age-column refresh.
$sheetNames = @("Standard Stats", "Shooting Stats", "Passing Stats", "Pass Types", "Goal & Shot Creation", "Defensive Actions", "Possession", "Playing Time", "Miscellaneous Stats")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 4; $r -le $lastRow; $r++) {
        $cellRef = "E" + $r
        $txt = $ws.Range($cellRef).Text
        if ($txt -match '^(\d+)-(\d+)$') {
            $yearPart = $Matches[1]
            $dayPart = [int]$Matches[2]
            $dayPart = $dayPart + 1
            if ($dayPart -lt 10) {
                $newVal = "$yearPart-00$dayPart"
            } elseif ($dayPart -lt 100) {
                $newVal = "$yearPart-0$dayPart"
            } else {
                $newVal = "$yearPart-$dayPart"
            }
            $ws.Range($cellRef).Value = $newVal
        }
    }
}
